# "add 2019 to files, update calibration"
#
# 1. EGGRA-gdp-adjustment: insert a new year column (2019) before the
#    existing 2020 column, fill header + value, autofit the new column.
# 2. About: the "Sources" cell (B4) loses its font/fill style override.
# 3. EGGRA-use-adjustment: the boolean "Use adjustment?" value flips from
#    1 to 0 (recalibration), and this sheet becomes the active tab with
#    a pending selection on B3 - mirroring EGGRA-gdp-adjustment's own
#    pending selection on B3.

$wb = $excel.ActiveWorkbook

# --- About -----------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("B4").Style = "Normal"

# --- EGGRA-gdp-adjustment ---------------------------------------------
$gdp = $wb.Worksheets.Item("EGGRA-gdp-adjustment")
$gdp.Columns("B").Insert()
$gdp.Range("B1").Value = 2019
$gdp.Range("B2").Value = 0
$gdp.Columns("B").AutoFit() | Out-Null
$gdp.Range("B3").Select() | Out-Null

# --- EGGRA-use-adjustment --------------------------------------------
$useAdj = $wb.Worksheets.Item("EGGRA-use-adjustment")
$useAdj.Range("B2").Value = 0
$useAdj.Activate() | Out-Null
$useAdj.Range("B3").Select() | Out-Null
